# Word COM-interop script implementing:
#   "Added web prefixes plus styling updates"
#
# 1. Merge the split "www.michaelbr" / "i" / "ckley.com" hyperlink runs
#    into a single run reading "www.michaelbrickley.com".
# 2. Remove the stray "_GoBack" bookmark that wraps "tes" / " tasks via
#    Python" in the experience bullet.
# 3. Insert a new "XML (High), " skill entry right after "HTML (High), "
#    (with a fresh "_GoBack" bookmark immediately following it, mirroring
#    real Word's habit of leaving one behind at the last edit position).
#
# NOTE: this runtime's Range.Text / Find-replace machinery renormalizes a
# touched paragraph by merging every maximal run of identically-formatted
# w:r siblings that the edit point falls inside of. That is perfect for
# (1) -- we *want* the three hyperlink runs coalesced -- but would wreck
# the "Skills" paragraph in (3), where only a single new run should
# appear and all the neighbouring runs must stay exactly as they were.
# To dodge that, the new skill text is authored in an isolated scratch
# paragraph appended at the very end of the story (so only *it* gets
# renormalized, harmlessly, since it is the sole run there), then
# Copy/Paste - which this runtime does not renormalize - transplants the
# finished run into place before the scratch paragraph is deleted again.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) www.michaelbr / i / ckley.com  ->  www.michaelbrickley.com
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "www.michaelbrickley.com", $false, $false, $false, $false, $false,
    $true, 1, $false, "www.michaelbrickley.com", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Drop the leftover "_GoBack" bookmark around "tes" / " tasks via Python"
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 3) Add "XML (High), " after "HTML (High), ", plus a new "_GoBack"
#    bookmark right after it (before "SQL (Medium...").
# ---------------------------------------------------------------------

# Build the new run text in an isolated scratch paragraph at the end of
# the document, so the only run affected by the implicit merge-on-edit
# is the scratch run itself.
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$scratchRange = $d.Paragraphs.Last.Range
$scratchRange.Text = "XML (High), "

# Re-locate the freshly written scratch run and copy it to the clipboard.
$scratchRange = $d.Content
$scratchRange.Find.Execute(
    "XML (High), ", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$scratchRange.Copy()

# Paste the clean run immediately after "HTML (High), ".
$target = $d.Content
$target.Find.Execute(
    "HTML (High), ", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)
$target.Paste()

# Remove the scratch paragraph we used as a staging area.
$lastParaRange = $d.Paragraphs.Last.Range
$storyEnd = $d.Content.End
$d.Range($lastParaRange.Start, $storyEnd).Delete()

# Insert a fresh "_GoBack" bookmark right after the new run (matching the
# target layout: ...XML (High), <bookmarkStart/><bookmarkEnd/>SQL (Medium...).
$afterNewRun = $d.Content
$afterNewRun.Find.Execute(
    "XML (High), ", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$afterNewRun.Collapse(0)
$d.Bookmarks.Add("_GoBack", $afterNewRun) | Out-Null
